$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(43, 3).Value = "passed"
$ws.Cells.Item(43, 4).Value = "20201116_204157"
$ws.Cells.Item(43, 5).Value = "chrome"

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(44, 3).Value = "failed"
$ws.Cells.Item(44, 4).Value = "20201116_204355"
$ws.Cells.Item(44, 5).Value = "chrome"

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(45, 3).Value = "passed"
$ws.Cells.Item(45, 4).Value = "20201116_204421"
$ws.Cells.Item(45, 5).Value = "chrome"

$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(46, 3).Value = "failed"
$ws.Cells.Item(46, 4).Value = "20201116_204544"
$ws.Cells.Item(46, 5).Value = "chrome"

$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(47, 3).Value = "failed"
$ws.Cells.Item(47, 4).Value = "20201116_210252"
$ws.Cells.Item(47, 5).Value = "chrome"

$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(48, 3).Value = "failed"
$ws.Cells.Item(48, 4).Value = "20201116_210352"
$ws.Cells.Item(48, 5).Value = "chrome"

$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(49, 3).Value = "passed"
$ws.Cells.Item(49, 4).Value = "20201116_210520"
$ws.Cells.Item(49, 5).Value = "chrome"

$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(50, 3).Value = "failed"
$ws.Cells.Item(50, 4).Value = "20201116_211551"
$ws.Cells.Item(50, 5).Value = "chrome"

$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(51, 3).Value = "failed"
$ws.Cells.Item(51, 4).Value = "20201116_211626"
$ws.Cells.Item(51, 5).Value = "chrome"

$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(52, 3).Value = "failed"
$ws.Cells.Item(52, 4).Value = "20201116_211702"
$ws.Cells.Item(52, 5).Value = "chrome"

$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(53, 3).Value = "failed"
$ws.Cells.Item(53, 4).Value = "20201116_211754"
$ws.Cells.Item(53, 5).Value = "chrome"

$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(54, 3).Value = "failed"
$ws.Cells.Item(54, 4).Value = "20201116_211954"
$ws.Cells.Item(54, 5).Value = "chrome"

$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(55, 3).Value = "failed"
$ws.Cells.Item(55, 4).Value = "20201116_212020"
$ws.Cells.Item(55, 5).Value = "chrome"

$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(56, 3).Value = "failed"
$ws.Cells.Item(56, 4).Value = "20201116_214112"
$ws.Cells.Item(56, 5).Value = "chrome"

$ws.Cells.Item(57, 1).Value = 56
$ws.Cells.Item(57, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(57, 3).Value = "failed"
$ws.Cells.Item(57, 4).Value = "20201116_220108"
$ws.Cells.Item(57, 5).Value = "chrome"

$ws.Cells.Item(58, 1).Value = 57
$ws.Cells.Item(58, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(58, 3).Value = "failed"
$ws.Cells.Item(58, 4).Value = "20201116_221813"
$ws.Cells.Item(58, 5).Value = "chrome"

$ws.Cells.Item(59, 1).Value = 58
$ws.Cells.Item(59, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(59, 3).Value = "failed"
$ws.Cells.Item(59, 4).Value = "20201116_221854"
$ws.Cells.Item(59, 5).Value = "chrome"

$ws.Cells.Item(60, 1).Value = 59
$ws.Cells.Item(60, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(60, 3).Value = "passed"
$ws.Cells.Item(60, 4).Value = "20201116_221931"
$ws.Cells.Item(60, 5).Value = "chrome"

$ws.Cells.Item(61, 1).Value = 60
$ws.Cells.Item(61, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(61, 3).Value = "failed"
$ws.Cells.Item(61, 4).Value = "20201116_222108"
$ws.Cells.Item(61, 5).Value = "chrome"

$ws.Cells.Item(62, 1).Value = 61
$ws.Cells.Item(62, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(62, 3).Value = "passed"
$ws.Cells.Item(62, 4).Value = "20201116_222350"
$ws.Cells.Item(62, 5).Value = "chrome"

$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(63, 3).Value = "failed"
$ws.Cells.Item(63, 4).Value = "20201118_220859"
$ws.Cells.Item(63, 5).Value = "chrome"

$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(64, 3).Value = "failed"
$ws.Cells.Item(64, 4).Value = "20201118_220952"
$ws.Cells.Item(64, 5).Value = "chrome"

$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(65, 3).Value = "failed"
$ws.Cells.Item(65, 4).Value = "20201118_222404"
$ws.Cells.Item(65, 5).Value = "chrome"

$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(66, 3).Value = "failed"
$ws.Cells.Item(66, 4).Value = "20201118_223525"
$ws.Cells.Item(66, 5).Value = "chrome"

$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(67, 3).Value = "failed"
$ws.Cells.Item(67, 4).Value = "20201118_223750"
$ws.Cells.Item(67, 5).Value = "chrome"

$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(68, 3).Value = "failed"
$ws.Cells.Item(68, 4).Value = "20201118_223909"
$ws.Cells.Item(68, 5).Value = "chrome"

$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(69, 3).Value = "passed"
$ws.Cells.Item(69, 4).Value = "20201118_224259"
$ws.Cells.Item(69, 5).Value = "chrome"

$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "testing-jdbc-country-functionality;country-tetsing-with-jdbc"
$ws.Cells.Item(70, 3).Value = "failed"
$ws.Cells.Item(70, 4).Value = "20201118_224839"
$ws.Cells.Item(70, 5).Value = "chrome"

$ws.Cells.Item(71, 1).Value = 70
$ws.Cells.Item(71, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(71, 3).Value = "failed"
$ws.Cells.Item(71, 4).Value = "20201118_224935"
$ws.Cells.Item(71, 5).Value = "chrome"

$ws.Cells.Item(72, 1).Value = 71
$ws.Cells.Item(72, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(72, 3).Value = "failed"
$ws.Cells.Item(72, 4).Value = "20201118_224956"
$ws.Cells.Item(72, 5).Value = "chrome"

$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(73, 3).Value = "failed"
$ws.Cells.Item(73, 4).Value = "20201118_225024"
$ws.Cells.Item(73, 5).Value = "chrome"

$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(74, 3).Value = "passed"
$ws.Cells.Item(74, 4).Value = "20201119_004141"
$ws.Cells.Item(74, 5).Value = "chrome"

$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(75, 3).Value = "passed"
$ws.Cells.Item(75, 4).Value = "20201119_004235"
$ws.Cells.Item(75, 5).Value = "chrome"

$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(76, 3).Value = "passed"
$ws.Cells.Item(76, 4).Value = "20201119_004610"
$ws.Cells.Item(76, 5).Value = "chrome"

$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(77, 3).Value = "passed"
$ws.Cells.Item(77, 4).Value = "20201119_004714"
$ws.Cells.Item(77, 5).Value = "chrome"

$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(78, 3).Value = "passed"
$ws.Cells.Item(78, 4).Value = "20201119_022749"
$ws.Cells.Item(78, 5).Value = "chrome"

$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(79, 3).Value = "passed"
$ws.Cells.Item(79, 4).Value = "20201119_023158"
$ws.Cells.Item(79, 5).Value = "chrome"

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(80, 3).Value = "passed"
$ws.Cells.Item(80, 4).Value = "20201119_023224"
$ws.Cells.Item(80, 5).Value = "chrome"

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(81, 3).Value = "passed"
$ws.Cells.Item(81, 4).Value = "20201119_023312"
$ws.Cells.Item(81, 5).Value = "chrome"

$ws.Cells.Item(82, 1).Value = 81
$ws.Cells.Item(82, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(82, 3).Value = "passed"
$ws.Cells.Item(82, 4).Value = "20201119_023621"
$ws.Cells.Item(82, 5).Value = "chrome"

$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(83, 3).Value = "passed"
$ws.Cells.Item(83, 4).Value = "20201119_023643"
$ws.Cells.Item(83, 5).Value = "chrome"

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(84, 3).Value = "failed"
$ws.Cells.Item(84, 4).Value = "20201119_023724"
$ws.Cells.Item(84, 5).Value = "chrome"

$ws.Cells.Item(85, 1).Value = 84
$ws.Cells.Item(85, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(85, 3).Value = "passed"
$ws.Cells.Item(85, 4).Value = "20201119_023750"
$ws.Cells.Item(85, 5).Value = "chrome"

$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = "testing-jdbc-city-functionality;city-tetsing-with-jdbc"
$ws.Cells.Item(86, 3).Value = "passed"
$ws.Cells.Item(86, 4).Value = "20201119_024223"
$ws.Cells.Item(86, 5).Value = "chrome"

$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "login-functionality;login-with-valid-username-and-password"
$ws.Cells.Item(87, 3).Value = "passed"
$ws.Cells.Item(87, 4).Value = "20201222_194026"
$ws.Cells.Item(87, 5).Value = "chrome"
